$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(44382,13101,69,6032,19202,5741,153,17,136,138),
    @(44383,13148,101,6055,19304,5768,149,17,132,138),
    @(44384,13251,101,6077,19429,5794,145,15,130,138),
    @(44385,13331,112,6104,19547,5800,166,15,151,138),
    @(44386,13423,84,6123,19630,5826,159,15,144,138),
    @(44387,13492,55,6145,19692,5846,161,12,149,138),
    @(44388,13510,40,6152,19702,5856,158,12,146,138),
    @(44389,13543,88,6163,19794,5867,158,9,149,138),
    @(44390,13593,136,6224,19953,5884,200,9,191,140),
    @(44391,13677,122,6252,20051,5890,221,9,212,141),
    @(44392,13789,112,6280,20181,5916,222,9,213,142),
    @(44393,13811,140,6289,20240,5939,208,11,197,142),
    @(44394,13929,89,6312,20330,5964,206,11,195,142),
    @(44395,13957,54,6319,20330,5984,193,11,182,142),
    @(44396,13972,84,6332,20388,6036,154,12,142,142),
    @(44397,14024,116,6348,20488,6066,140,12,128,142),
    @(44398,14099,115,6370,20584,6096,132,12,120,142),
    @(44399,14165,81,6386,20632,6096,148,10,138,142),
    @(44400,14208,87,6402,20697,6137,122,10,112,143),
    @(44401,14294,85,6424,20803,6158,123,10,113,143),
    @(44402,14294,85,6424,20803,6166,114,10,104,144),
    @(44403,14304,123,6435,20862,6176,115,10,105,144),
    @(44404,14401,110,6452,20963,6183,125,8,117,144)
)

$startRow = 369
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Activate()
$ws.Range("B392").Select()
